$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet as a copy of "2021-Q4" (same layout),
#    placed immediately before "总计", then overwrite its data.
# ---------------------------------------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy($null, $q4)
$q1 = $wb.Worksheets.Item(6)
$q1.Name = "2022-Q1"

# Remember the default (unstyled) look of the data cells so we can restore
# it after temporarily switching to a text format (this keeps values that
# look like numbers - e.g. "22.01" or "516110" - stored as text, exactly
# like the rest of the sheet, without leaving a stray custom style behind).
$bStyle = $q1.Range("B2:B4").Style
$dgStyle = $q1.Range("D2:G4").Style

$q1.Range("B2:B4").NumberFormat = "@"
$q1.Range("D2:G4").NumberFormat = "@"

$q1.Cells.Item(2, 2).Value = "004854"
$q1.Cells.Item(2, 3).Value = "广发中证全指汽车指数A"
$q1.Cells.Item(2, 4).Value = "22.01"
$q1.Cells.Item(2, 5).Value = "94.43"
$q1.Cells.Item(2, 6).Value = "4.58"
$q1.Cells.Item(2, 7).Value = "1.0081"

$q1.Cells.Item(3, 2).Value = "004855"
$q1.Cells.Item(3, 3).Value = "广发中证全指汽车指数C"
$q1.Cells.Item(3, 4).Value = "6.11"
$q1.Cells.Item(3, 5).Value = "94.43"
$q1.Cells.Item(3, 6).Value = "4.58"
$q1.Cells.Item(3, 7).Value = "0.2798"

$q1.Cells.Item(4, 2).Value = "516110"
$q1.Cells.Item(4, 3).Value = "国泰中证800汽车与零部件ETF"
$q1.Cells.Item(4, 4).Value = "1.74"
$q1.Cells.Item(4, 5).Value = "97.85"
$q1.Cells.Item(4, 6).Value = "2.73"
$q1.Cells.Item(4, 7).Value = "0.0475"

$q1.Range("B2:B4").Style = $bStyle
$q1.Range("D2:G4").Style = $dgStyle

# Rows/columns A and H already contain the correct values after the copy
# (0/1/2 and 7/7/10 respectively), so nothing else to change there.

# ---------------------------------------------------------------------
# 2. Insert a new first data row into "总计" for 2022-Q1, pushing the
#    existing rows down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

# Clear the formatting the insert copied down from the header row.
$total.Range("B2:D2").ClearFormats()

# A2 keeps the bold/centered style used by the rest of column A.
$total.Cells.Item(3, 1).Copy($total.Cells.Item(2, 1))
$total.Cells.Item(2, 1).Value = 0

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 3
$total.Cells.Item(2, 4).Value = 1.34

# Column A is a plain 0-based row counter (not a formula), so renumber the
# rows that shifted down by one position.
$total.Cells.Item(3, 1).Value = 1
$total.Cells.Item(4, 1).Value = 2
$total.Cells.Item(5, 1).Value = 3
$total.Cells.Item(6, 1).Value = 4
$total.Cells.Item(7, 1).Value = 5
